$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.063.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.70%  '
$ws.Range("D3").Value = "'2.758.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.35%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = "'580.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").Value = "'158.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.01%  '
$ws.Range("D7").Value = "'0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.49%  '
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = "'2.756.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.67%  '
$ws.Range("E10").Value = '  +3.77%  '
$ws.Range("E11").Value = '  +3.54%  '
$ws.Range("E12").Value = '  +3.63%  '
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").Value = "'3.228.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.66%  '
$ws.Range("D15").Value = "'27.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.68%  '
$ws.Range("D16").Value = "'63.962.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.62%  '
$ws.Range("E17").Value = '  +7.74%  '
$ws.Range("D18").Value = "'2.753.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.54%  '
$ws.Range("D19").Value = "'12.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.30%  '
$ws.Range("E20").Value = '  +4.10%  '
$ws.Range("D21").Value = "'364.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.08%  '
$ws.Range("D22").Value = "'7.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").Value = "'0.543"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.93%  '
$ws.Range("D24").Value = "'0.993"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("D25").Value = "'66.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.92%  '
$ws.Range("E26").Value = '  +5.80%  '
$ws.Range("D27").Value = "'8.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("D29").Value = "'0.0₃0916"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.82%  '
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").Value = "'7.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.31%  '
$ws.Range("E32").Value = '  +15.22%  '
$ws.Range("D33").Value = "'173.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.82%  '
$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = "'20.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.21%  '
$ws.Range("D36").Value = "'4.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.61%  '
$ws.Range("E37").Value = '  +8.21%  '
$ws.Range("E38").Value = '  +7.34%  '
$ws.Range("D39").Value = "'1.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.08%  '
$ws.Range("D40").Value = "'4.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("D41").Value = "'339.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = "'6.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +16.35%  '
$ws.Range("D43").Value = "'39.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.81%  '
$ws.Range("D44").Value = "'22.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.80%  '
$ws.Range("D45").Value = "'21.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.29%  '
$ws.Range("D46").Value = "'0.0606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.29%  '
$ws.Range("D47").Value = "'0.646"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("E48").Value = '  +3.17%  '
$ws.Range("D49").Value = "'137.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("D50").Value = "'0.102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.52%  '
$ws.Range("E51").Value = '  +0.28%  '
